$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values that look numeric stay as text (matches original inlineStr cells)
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "22.436.99"
$ws.Range("E2").Value = "  +0.15%  "

Set-TextValue $ws.Range("D3") "1.571.28"
$ws.Range("E3").Value = "  +0.59%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("E5").Value = "  +0.01%  "

Set-TextValue $ws.Range("D6") "287.62"
$ws.Range("E6").Value = "  +0.95%  "

Set-TextValue $ws.Range("D7") "0.3696"
$ws.Range("E7").Value = "  +2.01%  "

Set-TextValue $ws.Range("D8") "47.30"
$ws.Range("E8").Value = "  -2.21%  "

Set-TextValue $ws.Range("D9") "0.3313"
$ws.Range("E9").Value = "  -0.55%  "

$ws.Range("E10").Value = "  +2.60%  "

Set-TextValue $ws.Range("D11") "0.07497"
$ws.Range("E11").Value = "  +1.33%  "

$ws.Range("E12").Value = "  +0.08%  "

Set-TextValue $ws.Range("D13") "20.74"
$ws.Range("E13").Value = "  -0.22%  "

Set-TextValue $ws.Range("D14") "5.928"
$ws.Range("E14").Value = "  -0.01%  "

Set-TextValue $ws.Range("D15") "6.904"
$ws.Range("E15").Value = "  +0.24%  "

Set-TextValue $ws.Range("D16") "1.561.79"
$ws.Range("E16").Value = "  -0.10%  "

Set-TextValue $ws.Range("D17") "0.00001113"
$ws.Range("E17").Value = "  +0.89%  "

Set-TextValue $ws.Range("D18") "88.16"
$ws.Range("E18").Value = "  +0.06%  "

Set-TextValue $ws.Range("D19") "0.06723"
$ws.Range("E19").Value = "  +0.54%  "

Set-TextValue $ws.Range("D20") "6.406"
$ws.Range("E20").Value = "  +0.90%  "

$ws.Range("E21").Value = "  -0.03%  "

Set-TextValue $ws.Range("D22") "16.47"
$ws.Range("E22").Value = "  +2.27%  "

$ws.Range("E23").Value = "  -0.28%  "

Set-TextValue $ws.Range("D24") "22.419.61"
$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("E25").Value = "  -1.89%  "

Set-TextValue $ws.Range("D26") "2.621"
$ws.Range("E26").Value = "  +2.78%  "

Set-TextValue $ws.Range("D27") "150.43"
$ws.Range("E27").Value = "  +0.41%  "

$ws.Range("E28").Value = "  +0.86%  "

Set-TextValue $ws.Range("D29") "4.944"
$ws.Range("E29").Value = "  -0.99%  "

Set-TextValue $ws.Range("D30") "124.92"
$ws.Range("E30").Value = "  +1.55%  "

Set-TextValue $ws.Range("D31") "1.740.07"
$ws.Range("E31").Value = "  +0.13%  "

Set-TextValue $ws.Range("D32") "1.088"
$ws.Range("E32").Value = "  +2.75%  "

Set-TextValue $ws.Range("D35") "9.883"
$ws.Range("E35").Value = "  +0.80%  "

Set-TextValue $ws.Range("D36") "0.08318"
$ws.Range("E36").Value = "  +0.68%  "

Set-TextValue $ws.Range("D37") "0.02436"
$ws.Range("E37").Value = "  +1.65%  "

Set-TextValue $ws.Range("D38") "1.299"
$ws.Range("E38").Value = "  -0.28%  "

Set-TextValue $ws.Range("D39") "0.06369"
$ws.Range("E39").Value = "  -0.23%  "

Set-TextValue $ws.Range("D40") "0.2208"
$ws.Range("E40").Value = "  -0.05%  "

Set-TextValue $ws.Range("D41") "5.321"
$ws.Range("E41").Value = "  -0.02%  "

Set-TextValue $ws.Range("D42") "11.37"
$ws.Range("E42").Value = "  +2.18%  "

Set-TextValue $ws.Range("D43") "0.6220"
$ws.Range("E43").Value = "  +2.40%  "

$ws.Range("E44").Value = "  +0.04%  "

Set-TextValue $ws.Range("D45") "13.94"
$ws.Range("E45").Value = "  +0.71%  "

Set-TextValue $ws.Range("D46") "0.6032"
$ws.Range("E46").Value = "  +4.92%  "

Set-TextValue $ws.Range("D47") "3.770"
$ws.Range("E47").Value = "  +0.44%  "

$ws.Range("E48").Value = "  +1.30%  "

Set-TextValue $ws.Range("D49") "124.81"
$ws.Range("E49").Value = "  +0.20%  "

Set-TextValue $ws.Range("D50") "1.190"
$ws.Range("E50").Value = "  -1.62%  "

Set-TextValue $ws.Range("D51") "0.07192"
$ws.Range("E51").Value = "  -0.12%  "

# Rows 33 and 34 swap content (Filecoin <-> WEMIXTOKEN) plus value changes
$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D33") "1.988"
$ws.Range("E33").Value = "  -0.48%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D34") "6.070"
$ws.Range("E34").Value = "  -0.87%  "
